$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New reference materials data to append below the existing table (rows 16-22):
#   AlB2  (rows 16-19)
#   CaCu5 (rows 20-22)
$newRows = @(
  @("AlB2",  0.671, 8,  -8,  2),
  @("AlB2",  0.633, 13, -14, -1),
  @("AlB2",  0.58,  16, -20, -3),
  @("AlB2",  0.535, 20, -38, -17),
  @("CaCu5", 0.671, 6,  -8,  -1),
  @("CaCu5", 0.633, 10, -11, 0),
  @("CaCu5", 0.58,  11, -4,  6)
)

$r = 16
foreach ($row in $newRows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $r++
}

# Leave the selection on the last entered cell, like the source edit did.
$ws.Range("C22").Select() | Out-Null
